$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'316.35"
$ws.Range("E2").Value = "'1.51%"

$ws.Range("D3").Value = "'37.94"
$ws.Range("E3").Value = "'1.56%"

$ws.Range("D4").Value = "'5.173"
$ws.Range("E4").Value = "'0.71%"

$ws.Range("D5").Value = "'0.07962"
$ws.Range("E5").Value = "'1.78%"

$ws.Range("D6").Value = "'8.471"
$ws.Range("E6").Value = "'2.50%"

$ws.Range("D7").Value = "'1.920"
$ws.Range("E7").Value = "'0.97%"

$ws.Range("E8").Value = "'5.08%"

$ws.Range("D9").Value = "'0.9433"
$ws.Range("E9").Value = "'2.51%"

$ws.Range("D10").Value = "'0.1272"
$ws.Range("E10").Value = "'6.59%"

$ws.Range("D11").Value = "'0.1935"
$ws.Range("E11").Value = "'0.55%"

$ws.Range("D12").Value = "'0.08970"
$ws.Range("E12").Value = "'-1.73%"

$ws.Range("D13").Value = "'0.03433"
$ws.Range("E13").Value = "'2.42%"

$ws.Range("D14").Value = "'0.09533"
$ws.Range("E14").Value = "'-0.86%"

$ws.Range("D15").Value = "'0.001389"
$ws.Range("E15").Value = "'0.27%"

$ws.Range("D16").Value = "'0.006112"
$ws.Range("E16").Value = "'6.80%"

$ws.Range("D17").Value = "'3.421"
$ws.Range("E17").Value = "'-3.13%"

$ws.Range("D18").Value = "'4.469"
$ws.Range("E18").Value = "'1.17%"

$ws.Range("D19").Value = "'0.3513"
$ws.Range("E19").Value = "'2.08%"

$ws.Range("D20").Value = "'6.576"
$ws.Range("E20").Value = "'24.81%"

$ws.Range("E21").Value = "'1.65%"

$ws.Range("E23").Value = "'-0.04%"

$ws.Range("D24").Value = "'0.001218"
$ws.Range("E24").Value = "'-2.80%"

$ws.Range("D25").Value = "'0.004410"
$ws.Range("E25").Value = "'-5.52%"

$ws.Range("D26").Value = "'0.0001327"
$ws.Range("E26").Value = "'-2.75%"

$ws.Range("D27").Value = "'0.0003975"
$ws.Range("E27").Value = "'-0.70%"

$ws.Range("D39").Value = "'0.02399"
$ws.Range("E39").Value = "'5.39%"

$ws.Range("D40").Value = "'0.05163"
$ws.Range("E40").Value = "'2.09%"

$ws.Range("D41").Value = "'0.007447"
$ws.Range("E41").Value = "'-0.43%"

$ws.Range("D42").Value = "'0.1395"
$ws.Range("E42").Value = "'3.24%"

$ws.Range("D43").Value = "'0.008435"
$ws.Range("E43").Value = "'-7.24%"

$ws.Range("D44").Value = "'0.002115"
$ws.Range("E44").Value = "'8.10%"

$ws.Range("D45").Value = "'0.008753"
$ws.Range("E45").Value = "'-6.42%"

$ws.Range("D46").Value = "'0.00006494"
$ws.Range("E46").Value = "'-2.12%"

$ws.Range("D47").Value = "'0.00000000748"
$ws.Range("E47").Value = "'-0.56%"

$ws.Range("D48").Value = "'0.002861"
$ws.Range("E48").Value = "'-12.94%"

$ws.Range("D49").Value = "'0.001683"
$ws.Range("E49").Value = "'67.76%"

$ws.Range("D50").Value = "'0.00002095"
$ws.Range("E50").Value = "'-0.56%"

$ws.Range("D51").Value = "'0.0001996"
$ws.Range("E51").Value = "'-0.56%"
